$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("F2").Value = 20.81
$ws.Range("G2").Value = 14.2
$ws.Range("H2").Value = 83.66666666666667

# Row 3
$ws.Range("F3").Value = 54.79833333333332
$ws.Range("G3").Value = 8.699999999999999
$ws.Range("H3").Value = 63.66666666666666

# Row 4
$ws.Range("F4").Value = 58.04666666666667
$ws.Range("H4").Value = 68
